$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp shown in F1
$ws.Range("F1").Value = "Last status check on: 24.02.2022 09:45"

# Row 6 (Shell Olomoucká) got a fresh price reading:
#   - the new price goes into B6
#   - the previous B6 price shifts into C6 (the "old" price column)
#   - the delta column (D6) now stores a pre-formatted signed text value
#   - the date column (E6) now stores a plain text timestamp instead of a date serial
$ws.Range("B6").Value = 38.29
$ws.Range("C6").Value = 37.9

$d6 = $ws.Range("D6")
$d6.NumberFormat = "@"
$d6.Value = "+0.39"

$e6 = $ws.Range("E6")
$e6.NumberFormat = "@"
$e6.Value = "2022-02-24 09:47:32"
